$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 1929
$ws.Range("C4").Value = 80.40000000000001
$ws.Range("D4").Value = 336.4
$ws.Range("I4").Value = 1648
$ws.Range("K4").Value = 482.2
$ws.Range("N4").Value = -1.6
$ws.Range("O4").Value = 698
$ws.Range("P4").Value = 29.1
$ws.Range("Q4").Value = 131
$ws.Range("V4").Value = 643
$ws.Range("X4").Value = 174.5
$ws.Range("AB4").Value = 2285
$ws.Range("AC4").Value = 95.2
$ws.Range("AD4").Value = 250.7
$ws.Range("AH4").Value = 7.2
$ws.Range("AI4").Value = 981
$ws.Range("AK4").Value = 380.8
$ws.Range("AL4").Value = 6
$ws.Range("AM4").Value = 25
$ws.Range("AN4").Value = -1.5
$ws.Range("B5").Value = 2217968
$ws.Range("C5").Value = 885.4
$ws.Range("D5").Value = 2152.4
$ws.Range("G5").Value = 132
$ws.Range("H5").Value = 1024
$ws.Range("I5").Value = 55081
$ws.Range("K5").Value = 1413.6
$ws.Range("L5").Value = 1569
$ws.Range("M5").Value = 62.6
$ws.Range("O5").Value = 860841
$ws.Range("P5").Value = 343.6
$ws.Range("Q5").Value = 1581
$ws.Range("T5").Value = 19
$ws.Range("U5").Value = 190
$ws.Range("V5").Value = 48717
$ws.Range("X5").Value = 537
$ws.Range("Y5").Value = 1603
$ws.Range("Z5").Value = 64
$ws.Range("AA5").Value = 1.2
$ws.Range("AB5").Value = 1136889
$ws.Range("AC5").Value = 453.8
$ws.Range("AD5").Value = 645.5
$ws.Range("AG5").Value = 175
$ws.Range("AH5").Value = 706
$ws.Range("AI5").Value = 5487
$ws.Range("AK5").Value = 685.7
$ws.Range("AL5").Value = 1658
$ws.Range("AM5").Value = 66.2
$ws.Range("AN5").Value = 0.9
$ws.Range("B6").Value = 55848
$ws.Range("C6").Value = 3723.2
$ws.Range("D6").Value = 13441.8
$ws.Range("I6").Value = 52263
$ws.Range("K6").Value = 11169.6
$ws.Range("O6").Value = 25527
$ws.Range("P6").Value = 1701.8
$ws.Range("Q6").Value = 6341.5
$ws.Range("U6").Value = 10.5
$ws.Range("V6").Value = 24618
$ws.Range("X6").Value = 5105.4
$ws.Range("AA6").Value = -0.6
$ws.Range("AB6").Value = 1993
$ws.Range("AC6").Value = 132.9
$ws.Range("AD6").Value = 305.3
$ws.Range("AH6").Value = 56
$ws.Range("AI6").Value = 932
$ws.Range("AK6").Value = 398.6
$ws.Range("AL6").Value = 5
$ws.Range("AM6").Value = 33.3
$ws.Range("AN6").Value = -1
$ws.Range("B7").Value = 66181
$ws.Range("C7").Value = 704.1
$ws.Range("D7").Value = 2557
$ws.Range("H7").Value = 567.8
$ws.Range("I7").Value = 22972
$ws.Range("K7").Value = 1504.1
$ws.Range("L7").Value = 44
$ws.Range("M7").Value = 46.8
$ws.Range("O7").Value = 26802
$ws.Range("P7").Value = 285.1
$ws.Range("Q7").Value = 1828.3
$ws.Range("U7").Value = 84.5
$ws.Range("V7").Value = 17700
$ws.Range("X7").Value = 638.1
$ws.Range("Y7").Value = 42
$ws.Range("Z7").Value = 44.7
$ws.Range("AA7").Value = 0
$ws.Range("AB7").Value = 34350
$ws.Range("AC7").Value = 365.4
$ws.Range("AD7").Value = 451
$ws.Range("AG7").Value = 175
$ws.Range("AH7").Value = 701.5
$ws.Range("AI7").Value = 1627
$ws.Range("AK7").Value = 648.1
$ws.Range("AL7").Value = 53
$ws.Range("AM7").Value = 56.4
$ws.Range("AN7").Value = 0.4
$ws.Range("B8").Value = 365975
$ws.Range("C8").Value = 740.8
$ws.Range("D8").Value = 2684.7
$ws.Range("G8").Value = 0.5
$ws.Range("H8").Value = 625.5
$ws.Range("I8").Value = 50844
$ws.Range("K8").Value = 1481.7
$ws.Range("L8").Value = 247
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 0.3
$ws.Range("O8").Value = 117975
$ws.Range("P8").Value = 238.8
$ws.Range("Q8").Value = 1021.5
$ws.Range("T8").Value = 1.5
$ws.Range("U8").Value = 107
$ws.Range("V8").Value = 11140
$ws.Range("X8").Value = 466.3
$ws.Range("Y8").Value = 253
$ws.Range("Z8").Value = 51.2
$ws.Range("AB8").Value = 197374
$ws.Range("AC8").Value = 399.5
$ws.Range("AD8").Value = 743
$ws.Range("AG8").Value = 72
$ws.Range("AH8").Value = 561
$ws.Range("AI8").Value = 8295
$ws.Range("AK8").Value = 675.9
$ws.Range("AL8").Value = 292
$ws.Range("AM8").Value = 59.1
$ws.Range("AN8").Value = 0.5
$ws.Range("B9").Value = 115047
$ws.Range("C9").Value = 555.8
$ws.Range("D9").Value = 1107.7
$ws.Range("H9").Value = 597.5
$ws.Range("I9").Value = 6588
$ws.Range("K9").Value = 991.8
$ws.Range("O9").Value = 22411
$ws.Range("P9").Value = 108.3
$ws.Range("Q9").Value = 346.9
$ws.Range("T9").Value = 4
$ws.Range("U9").Value = 82.5
$ws.Range("V9").Value = 4083
$ws.Range("X9").Value = 193.2
$ws.Range("Y9").Value = 116
$ws.Range("Z9").Value = 56
$ws.Range("AB9").Value = 91795
$ws.Range("AC9").Value = 443.5
$ws.Range("AD9").Value = 675.4
$ws.Range("AG9").Value = 155
$ws.Range("AH9").Value = 627
$ws.Range("AI9").Value = 3597
$ws.Range("AK9").Value = 717.1
$ws.Range("AL9").Value = 128
$ws.Range("AM9").Value = 61.8
$ws.Range("AN9").Value = 0.7
